{"js": "// The document had two small text corrections made to it:\n//   1. \"Xen Virtualization\" -> \"Zen Virtualization\"   (typo fix in the\n//      \"Virtualization:\" skills line)\n//   2. \"2019\" -> \"2009\"   (graduation year for \"The Swaminarayan School,\n//      NAGPUR\")\n//\n// (A full re-save in Word also dropped a number of now-stale\n// <w:proofErr> spell-check markers sprinkled across the document and\n// re-flowed the runs around them, but that produces no visible text\n// change, so we only need to make the two real edits below.)\n\nconst body = context.document.body;\n\n// 1) \"Xen Virtualization\" -> \"Zen Virtualization\"\nconst xenResults = body.search(\"Xen Virtualization\", { matchCase: true, matchWholeWord: false });\nxenResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < xenResults.items.length; i++) {\n  xenResults.items[i].insertText(\"Zen Virtualization\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) \"2019\" -> \"2009\"\nconst yearResults = body.search(\"2019\", { matchCase: true, matchWholeWord: false });\nyearResults.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < yearResults.items.length; i++) {\n  yearResults.items[i].insertText(\"2009\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The document had two small text corrections made to it:\n#   1. \"Xen Virtualization\" -> \"Zen Virtualization\"   (typo fix in the\n#      \"Virtualization:\" skills line)\n#   2. \"2019\" -> \"2009\"   (graduation year for \"The Swaminarayan School,\n#      NAGPUR\")\n#\n# (A full re-save in Word also dropped a number of now-stale\n# <w:proofErr/> spell-check markers sprinkled across the document and\n# re-flowed the runs around them, but that produces no visible text\n# change, so we only need to make the two real edits below.)\n\n$d = $word.ActiveDocument\n\n# wdReplaceAll = 2\n$wdReplaceAll = 2\n\n# 1) \"Xen Virtualization\" -> \"Zen Virtualization\"\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"Xen Virtualization\"\n$find1.Replacement.Text = \"Zen Virtualization\"\n$find1.Execute($find1.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find1.Replacement.Text, $wdReplaceAll)\n\n# 2) \"2019\" -> \"2009\"\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"2019\"\n$find2.Replacement.Text = \"2009\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, $wdReplaceAll)\n"}
